# Editing the methods section.
# Mark the two comments about "material exergy" definition (rows 8 and 9 of
# "Review 1") as resolved: copy the "DONE" row formatting used elsewhere in
# the sheet, fill in the Response (F) / Status (G) columns, switch the
# responder (I) from ES to MKH, and update the row height / selection to
# match the edited review sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Review 1")

# --- Row 8 -----------------------------------------------------------
# Re-use the exact cell formatting already used by the other resolved
# ("DONE") rows (e.g. row 4) so the styles line up with the rest of the
# sheet instead of inventing new ones.
$ws.Range("A4:I4").Copy()
$ws.Range("A8:I8").PasteSpecial(-4122)

$ws.Range("F8").Value = "Yes. Laura is right. I made a small change to make this clearer"
$ws.Range("G8").Value = "DONE"
$ws.Range("I8").Value = "MKH"

# --- Row 9 -------------------------------------------------------------
$ws.Range("A4:I4").Copy()
$ws.Range("A9:I9").PasteSpecial(-4122)

$responseText9 = "Didn't quite make this change, but did something better. Material exergy is defined as the maximum amount of work`nthat could be extracted`nby an ideal, reversible process`nthat brings the mixture of materials`ninto equilibrium with the reference environment,`ncharacterized by`nparticle size (`$d_0`$),`ntemperature (`$T_0`$),`npressure (`$P_0`$), and`nchemical composition (`$y_{i,0}`$)."
$ws.Range("F9").Value = $responseText9
$ws.Range("G9").Value = "DONE"
$ws.Range("I9").Value = "MKH"

# Row 9 grew substantially once the long response was added & reformatted.
$ws.Rows(9).RowHeight = 176

# --- Selection -----------------------------------------------------------
$ws.Range("E10").Select()
